$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.ClearContents()

# Header row (row 1)
$ws.Range("A1").Value = 'Sending cluster'
$ws.Range("B1").Value = 'Ligand symbol'
$ws.Range("C1").Value = 'Receptor symbol'
$ws.Range("D1").Value = 'Target cluster'
$ws.Range("E1").Value = 'Ligand-expressing cells'
$ws.Range("F1").Value = 'Ligand detection rate'
$ws.Range("G1").Value = 'Ligand average expression value'
$ws.Range("H1").Value = 'Ligand total expression value'
$ws.Range("I1").Value = 'Ligand derived specificity of average expression value'
$ws.Range("J1").Value = 'Ligand derived specificity of total expression value'
$ws.Range("K1").Value = 'Receptor-expressing cells'
$ws.Range("L1").Value = 'Receptor detection rate'
$ws.Range("M1").Value = 'Receptor average expression value'
$ws.Range("N1").Value = 'Receptor total expression value'
$ws.Range("O1").Value = 'Receptor derived specificity of average expression value'
$ws.Range("P1").Value = 'Receptor derived specificity of total expression value'
$ws.Range("Q1").Value = 'Edge average expression weight'
$ws.Range("R1").Value = 'Edge total expression weight'
$ws.Range("S1").Value = 'Edge average expression derived specificity'
$ws.Range("T1").Value = 'Edge total expression derived specificity'

# Data columns A-D, written column-by-column so new shared strings are
# appended to the string table in the same order pandas/openpyxl produced
# them originally (column-major first-use order).
# Column A
$ws.Range("A2").Value = 'ECs'
$ws.Range("A3").Value = 'ECs'
$ws.Range("A4").Value = 'FAPs'
$ws.Range("A5").Value = 'FAPs'
$ws.Range("A6").Value = 'M2'
$ws.Range("A7").Value = 'M2'
$ws.Range("A8").Value = 'sCs'
$ws.Range("A9").Value = 'sCs'

# Column B
$ws.Range("B2").Value = 'Adm'
$ws.Range("B3").Value = 'Adm'
$ws.Range("B4").Value = 'Adm'
$ws.Range("B5").Value = 'Adm'
$ws.Range("B6").Value = 'Adm'
$ws.Range("B7").Value = 'Adm'
$ws.Range("B8").Value = 'Adm'
$ws.Range("B9").Value = 'Adm'

# Column C
$ws.Range("C2").Value = 'Ramp3'
$ws.Range("C3").Value = 'Ramp3'
$ws.Range("C4").Value = 'Ramp3'
$ws.Range("C5").Value = 'Ramp3'
$ws.Range("C6").Value = 'Ramp3'
$ws.Range("C7").Value = 'Ramp3'
$ws.Range("C8").Value = 'Ramp3'
$ws.Range("C9").Value = 'Ramp3'

# Column D
$ws.Range("D2").Value = 'ECs'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("D4").Value = 'ECs'
$ws.Range("D5").Value = 'FAPs'
$ws.Range("D6").Value = 'ECs'
$ws.Range("D7").Value = 'FAPs'
$ws.Range("D8").Value = 'ECs'
$ws.Range("D9").Value = 'FAPs'

# Remaining numeric columns E-T, row by row
# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.73555533333333
$ws.Range("H2").Value = 56.206666
$ws.Range("I2").Value = 0.4699290876663871
$ws.Range("J2").Value = 0.4699290876663871
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.240458
$ws.Range("N2").Value = 3.721374
$ws.Range("O2").Value = 0.7314599148753498
$ws.Range("P2").Value = 0.7314599148753498
$ws.Range("Q2").Value = 23.240669497676
$ws.Range("R2").Value = 209.166025479084
$ws.Range("S2").Value = 0.3437342904619063
$ws.Range("T2").Value = 0.3437342904619063

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.73555533333333
$ws.Range("H3").Value = 56.206666
$ws.Range("I3").Value = 0.4699290876663871
$ws.Range("J3").Value = 0.4699290876663871
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.455408
$ws.Range("N3").Value = 1.366224
$ws.Range("O3").Value = 0.2685400851246502
$ws.Range("P3").Value = 0.2685400851246502
$ws.Range("Q3").Value = 8.532321783242667
$ws.Range("R3").Value = 76.790896049184
$ws.Range("S3").Value = 0.1261947972044808
$ws.Range("T3").Value = 0.1261947972044808

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 19.55844
$ws.Range("H4").Value = 58.67532
$ws.Range("I4").Value = 0.490568851675588
$ws.Range("J4").Value = 0.4905688516755881
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.240458
$ws.Range("N4").Value = 3.721374
$ws.Range("O4").Value = 0.7314599148753498
$ws.Range("P4").Value = 0.7314599148753498
$ws.Range("Q4").Value = 24.26142336552
$ws.Range("R4").Value = 218.35281028968
$ws.Range("S4").Value = 0.3588314504871237
$ws.Range("T4").Value = 0.3588314504871237

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.55844
$ws.Range("H5").Value = 58.67532
$ws.Range("I5").Value = 0.490568851675588
$ws.Range("J5").Value = 0.4905688516755881
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.455408
$ws.Range("N5").Value = 1.366224
$ws.Range("O5").Value = 0.2685400851246502
$ws.Range("P5").Value = 0.2685400851246502
$ws.Range("Q5").Value = 8.907070043520001
$ws.Range("R5").Value = 80.16363039168
$ws.Range("S5").Value = 0.1317374011884643
$ws.Range("T5").Value = 0.1317374011884643

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3162143333333333
$ws.Range("H6").Value = 0.9486429999999999
$ws.Range("I6").Value = 0.007931353542853873
$ws.Range("J6").Value = 0.007931353542853875
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.240458
$ws.Range("N6").Value = 3.721374
$ws.Range("O6").Value = 0.7314599148753498
$ws.Range("P6").Value = 0.7314599148753498
$ws.Range("Q6").Value = 0.392250599498
$ws.Range("R6").Value = 3.530255395482
$ws.Range("S6").Value = 0.005801467187302198
$ws.Range("T6").Value = 0.005801467187302199

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3162143333333333
$ws.Range("H7").Value = 0.9486429999999999
$ws.Range("I7").Value = 0.007931353542853873
$ws.Range("J7").Value = 0.007931353542853875
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.455408
$ws.Range("N7").Value = 1.366224
$ws.Range("O7").Value = 0.2685400851246502
$ws.Range("P7").Value = 0.2685400851246502
$ws.Range("Q7").Value = 0.1440065371146667
$ws.Range("R7").Value = 1.296058834032
$ws.Range("S7").Value = 0.002129886355551675
$ws.Range("T7").Value = 0.002129886355551676

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 1.258689333333333
$ws.Range("H8").Value = 3.776068
$ws.Range("I8").Value = 0.03157070711517098
$ws.Range("J8").Value = 0.03157070711517099
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.240458
$ws.Range("N8").Value = 3.721374
$ws.Range("O8").Value = 0.7314599148753498
$ws.Range("P8").Value = 0.7314599148753498
$ws.Range("Q8").Value = 1.561351253048
$ws.Range("R8").Value = 14.052161277432
$ws.Range("S8").Value = 0.02309270673901756
$ws.Range("T8").Value = 0.02309270673901757

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 1.258689333333333
$ws.Range("H9").Value = 3.776068
$ws.Range("I9").Value = 0.03157070711517098
$ws.Range("J9").Value = 0.03157070711517099
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.455408
$ws.Range("N9").Value = 1.366224
$ws.Range("O9").Value = 0.2685400851246502
$ws.Range("P9").Value = 0.2685400851246502
$ws.Range("Q9").Value = 0.5732171919146667
$ws.Range("R9").Value = 5.158954727232
$ws.Range("S9").Value = 0.008478000376153413
$ws.Range("T9").Value = 0.008478000376153417
